# Inserta una nueva columna "MAE" entre "R2" y "Tipo"
# (grafico de tipos de modelo - preparacion de datos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insertar una nueva columna en la posicion D, desplazando la columna
# "Tipo" (y sus datos) una posicion a la derecha, hacia la columna E.
$ws.Columns.Item(4).Insert()

# Copiar el formato del encabezado existente (C1, estilo con negrita y
# borde) hacia la nueva celda de encabezado D1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Rellenar el nuevo encabezado y el valor de la metrica MAE.
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 1.863313287418611
